$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H31").Value = 900
$ws.Range("J31").Value = 1500
$ws.Range("L31").Value = 4500
$ws.Range("N31").Value = -4960
$ws.Range("H40").Value = 1250
$ws.Range("I40").Value = 1000
$ws.Range("J40").Value = 1500
$ws.Range("K40").Value = 1000
$ws.Range("L40").Value = 1500
$ws.Range("M40").Value = -825
$ws.Range("N40").Value = -1850
$ws.Range("H64").Value = 3120
$ws.Range("I64").Value = 3200
$ws.Range("J64").Value = 3000
$ws.Range("K64").Value = 3200
$ws.Range("L64").Value = 3000
$ws.Range("M64").Value = -2952
$ws.Range("N64").Value = -3496
$ws.Range("H67").Value = 3120
$ws.Range("I67").Value = 3200
$ws.Range("J67").Value = 3000
$ws.Range("K67").Value = 3200
$ws.Range("L67").Value = 3000
$ws.Range("M67").Value = -2342
$ws.Range("N67").Value = -4716
$ws.Range("H70").Value = 1539.8
$ws.Range("I70").Value = 1499.8572
$ws.Range("J70").Value = 1633
$ws.Range("K70").Value = 4499.571599999999
$ws.Range("L70").Value = 4899
$ws.Range("M70").Value = -4229.571599999999
$ws.Range("N70").Value = -5439
$ws.Range("H73").Value = 1539.8
$ws.Range("I73").Value = 1499.8572
$ws.Range("J73").Value = 1633
$ws.Range("K73").Value = 4499.571599999999
$ws.Range("L73").Value = 4899
$ws.Range("M73").Value = -3563.571599999999
$ws.Range("N73").Value = -6771
$ws.Range("H76").Value = 3363.6365
$ws.Range("I76").Value = 3363.6365
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = 3363.6365
$ws.Range("L76").Value = 0
$ws.Range("M76").Value = -3048.6365
$ws.Range("N76").Value = ""
$ws.Range("H79").Value = 3363.6365
$ws.Range("I79").Value = 3363.6365
$ws.Range("J79").Value = 0
$ws.Range("K79").Value = 3363.6365
$ws.Range("L79").Value = 0
$ws.Range("M79").Value = -2271.6365
$ws.Range("N79").Value = ""
$ws.Range("H82").Value = 1679.3334
$ws.Range("I82").Value = 1679.3334
$ws.Range("K82").Value = 5038.0002
$ws.Range("M82").Value = -4632.0002
$ws.Range("H85").Value = 1679.3334
$ws.Range("I85").Value = 1679.3334
$ws.Range("K85").Value = 5038.0002
$ws.Range("M85").Value = -3634.0002
$ws.Range("H88").Value = 16157.45
$ws.Range("I88").Value = 856.4286
$ws.Range("J88").Value = 24396.46
$ws.Range("K88").Value = 856.4286
$ws.Range("L88").Value = 24396.46
$ws.Range("M88").Value = -450.4286
$ws.Range("N88").Value = -25208.46
$ws.Range("H91").Value = 16157.45
$ws.Range("I91").Value = 856.4286
$ws.Range("J91").Value = 24396.46
$ws.Range("K91").Value = 856.4286
$ws.Range("L91").Value = 24396.46
$ws.Range("M91").Value = 547.5714
$ws.Range("N91").Value = -27204.46
$ws.Range("H115").Value = 3364
$ws.Range("I115").Value = 3728
$ws.Range("K115").Value = 11184
$ws.Range("M115").Value = -9617
$ws.Range("H125").Value = 1900.8889
$ws.Range("I125").Value = 2000
$ws.Range("J125").Value = 1851.3334
$ws.Range("K125").Value = 18000
$ws.Range("L125").Value = 16662.0006
$ws.Range("M125").Value = -15540
$ws.Range("N125").Value = -21582.0006
$ws.Range("H129").Value = 1141.7637
$ws.Range("I129").Value = 556.6667
$ws.Range("J129").Value = 1256.2391
$ws.Range("K129").Value = 1670.0001
$ws.Range("L129").Value = 3768.7173
$ws.Range("M129").Value = 3329.9999
$ws.Range("N129").Value = -13768.7173
$ws.Range("H138").Value = 2150.7812
$ws.Range("J138").Value = 2135.6365
$ws.Range("L138").Value = 6406.9095
$ws.Range("N138").Value = -16686.9095

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 6434.625
$ws.Range("I63").Value = 3000
$ws.Range("J63").Value = 7579.5
$ws.Range("K63").Value = 3000
$ws.Range("L63").Value = 7579.5
$ws.Range("M63").Value = -2314
$ws.Range("N63").Value = -8951.5
$ws.Range("H66").Value = 6434.625
$ws.Range("I66").Value = 3000
$ws.Range("J66").Value = 7579.5
$ws.Range("K66").Value = 15000
$ws.Range("L66").Value = 37897.5
$ws.Range("M66").Value = -11568
$ws.Range("N66").Value = -44761.5
$ws.Range("H88").Value = 3129.75
$ws.Range("J88").Value = 3007
$ws.Range("L88").Value = 3007
$ws.Range("N88").Value = -3819
$ws.Range("H91").Value = 3129.75
$ws.Range("J91").Value = 3007
$ws.Range("L91").Value = 3007
$ws.Range("N91").Value = -5815
$ws.Range("H102").Value = 1857.3684
$ws.Range("I102").Value = 1722.9412
$ws.Range("K102").Value = 1722.9412
$ws.Range("M102").Value = -100.9412

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H26").Value = 23700
$ws.Range("I26").Value = 23700
$ws.Range("K26").Value = 23700
$ws.Range("M26").Value = -23408
$ws.Range("H80").Value = 1234.2354
$ws.Range("J80").Value = 153.08333
$ws.Range("L80").Value = 153.08333
$ws.Range("N80").Value = -2149.08333
$ws.Range("H83").Value = 1234.2354
$ws.Range("J83").Value = 153.08333
$ws.Range("L83").Value = 765.4166499999999
$ws.Range("N83").Value = -10749.41665
$ws.Range("H86").Value = 58825612
$ws.Range("I86").Value = 58825612
$ws.Range("K86").Value = 58825612
$ws.Range("M86").Value = -58824489
$ws.Range("H89").Value = 58825612
$ws.Range("I89").Value = 58825612
$ws.Range("K89").Value = 294128060
$ws.Range("M89").Value = -294122444

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 8774141
$ws.Range("I132").Value = 1731.7693
$ws.Range("K132").Value = 5195.3079
$ws.Range("M132").Value = -2665.3079
$ws.Range("H133").Value = 19800
$ws.Range("J133").Value = 19800
$ws.Range("L133").Value = 19800
$ws.Range("N133").Value = -24860

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H110").Value = 10122.652
$ws.Range("J110").Value = 11094.421
$ws.Range("L110").Value = 33283.263
$ws.Range("N110").Value = -41463.263
$ws.Range("H113").Value = 937.7857
$ws.Range("J113").Value = 1052.6364
$ws.Range("L113").Value = 3157.9092
$ws.Range("N113").Value = -7497.9092
$ws.Range("H115").Value = 4958.5557
$ws.Range("I115").Value = 3916.25
$ws.Range("J115").Value = 5792.4
$ws.Range("K115").Value = 11748.75
$ws.Range("L115").Value = 17377.2
$ws.Range("M115").Value = -10573.75
$ws.Range("N115").Value = -19727.2
$ws.Range("H119").Value = 3175
$ws.Range("I119").Value = 2566.6667
$ws.Range("K119").Value = 7700.000100000001
$ws.Range("M119").Value = -2862.000100000001
$ws.Range("H137").Value = 5562540.5
$ws.Range("I137").Value = 15163031
$ws.Range("J137").Value = 4362
$ws.Range("K137").Value = 45489093
$ws.Range("L137").Value = 13086
$ws.Range("M137").Value = -45483993
$ws.Range("N137").Value = -23286
$ws.Range("H140").Value = 1497.5927
$ws.Range("I140").Value = 1232.3334
$ws.Range("J140").Value = 2028.1111
$ws.Range("K140").Value = 3697.0002
$ws.Range("L140").Value = 6084.3333
$ws.Range("M140").Value = 1482.9998
$ws.Range("N140").Value = -16444.3333

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6651.231
$ws.Range("I70").Value = 6666.6
$ws.Range("K70").Value = 6666.6
$ws.Range("M70").Value = -6396.6
$ws.Range("H73").Value = 6651.231
$ws.Range("I73").Value = 6666.6
$ws.Range("K73").Value = 6666.6
$ws.Range("M73").Value = -5730.6
$ws.Range("H80").Value = 34207.69
$ws.Range("I80").Value = 4450
$ws.Range("J80").Value = 59714.285
$ws.Range("K80").Value = 4450
$ws.Range("L80").Value = 59714.285
$ws.Range("M80").Value = -3452
$ws.Range("N80").Value = -61710.285
$ws.Range("H83").Value = 34207.69
$ws.Range("I83").Value = 4450
$ws.Range("J83").Value = 59714.285
$ws.Range("K83").Value = 22250
$ws.Range("L83").Value = 298571.425
$ws.Range("M83").Value = -17258
$ws.Range("N83").Value = -308555.425
$ws.Range("H126").Value = 1915.5714
$ws.Range("I126").Value = 1665.2142
$ws.Range("J126").Value = 2416.2856
$ws.Range("K126").Value = 4995.642599999999
$ws.Range("L126").Value = 7248.8568
$ws.Range("M126").Value = -2525.642599999999
$ws.Range("N126").Value = -12188.8568

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2184.8333
$ws.Range("I7").Value = 1576
$ws.Range("J7").Value = 3402.5
$ws.Range("K7").Value = 1576
$ws.Range("L7").Value = 3402.5
$ws.Range("M7").Value = -1464
$ws.Range("N7").Value = -3626.5
$ws.Range("H22").Value = 9684.190000000001
$ws.Range("I22").Value = 625.1667
$ws.Range("K22").Value = 625.1667
$ws.Range("M22").Value = -330.1667
$ws.Range("H27").Value = 9684.190000000001
$ws.Range("I27").Value = 625.1667
$ws.Range("K27").Value = 625.1667
$ws.Range("M27").Value = -518.1667
$ws.Range("H68").Value = 2285.0938
$ws.Range("I68").Value = 1898
$ws.Range("K68").Value = 1898
$ws.Range("M68").Value = -1149
$ws.Range("H71").Value = 2285.0938
$ws.Range("I71").Value = 1898
$ws.Range("K71").Value = 9490
$ws.Range("M71").Value = -5746
$ws.Range("H126").Value = 2184.8333
$ws.Range("I126").Value = 1576
$ws.Range("J126").Value = 3402.5
$ws.Range("K126").Value = 4728
$ws.Range("L126").Value = 10207.5
$ws.Range("M126").Value = -2258
$ws.Range("N126").Value = -15147.5
$ws.Range("H132").Value = 4493.8335
$ws.Range("I132").Value = 3456.5557
$ws.Range("J132").Value = 5531.1113
$ws.Range("K132").Value = 10369.6671
$ws.Range("L132").Value = 16593.3339
$ws.Range("M132").Value = -7839.667099999999
$ws.Range("N132").Value = -21653.3339

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 10002
$ws.Range("I100").Value = 10002
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 20004
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -19463
$ws.Range("N100").Value = ""
$ws.Range("H126").Value = 1529.1428
$ws.Range("I126").Value = 1529.1428
$ws.Range("K126").Value = 4587.428400000001
$ws.Range("M126").Value = -2117.428400000001
$ws.Range("H132").Value = 4388006.5
$ws.Range("I132").Value = 2118
$ws.Range("J132").Value = 6174850.5
$ws.Range("K132").Value = 6354
$ws.Range("L132").Value = 18524551.5
$ws.Range("M132").Value = -3824
$ws.Range("N132").Value = -18529611.5

Write-Output "Applied changes"